$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '43.769.74'
$ws.Range("E2").Value = '  +0.11%  '

# Row 3
$ws.Range("D3").Value = '2.289.98'
$ws.Range("E3").Value = '  -0.09%  '

# Row 4
$ws.Range("E4").Value = '  +0.19%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '115.46'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +16.91%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '269.15'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.57%  '

# Row 7
$ws.Range("E7").Value = '  +0.16%  '

# Row 8
$ws.Range("E8").Value = '  +0.06%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.618'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.20%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '48.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +8.47%  '

# Row 11
$ws.Range("E11").Value = '  +1.05%  '

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.04'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +14.40%  '

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.107'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.47%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '15.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.72%  '

# Row 15
$ws.Range("D15").Value = '2.634.63'
$ws.Range("E15").Value = '  -0.04%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.877'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.82%  '

# Row 17
$ws.Range("D17").Value = '2.286.63'
$ws.Range("E17").Value = '  -0.16%  '

# Row 18
$ws.Range("D18").Value = '43.655.66'
$ws.Range("E18").Value = '  -0.23%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.0000109'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.03'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +13.49%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '72.27'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.23%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.42'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.92%  '

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +10.09%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '233.13'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.04%  '

# Row 25
$ws.Range("E25").Value = '  +3.33%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '11.75'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.83%  '

# Row 27
$ws.Range("E27").Value = '  +0.01%  '

# Row 28
$ws.Range("B28").Value = 'LEO'
$ws.Range("C28").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.92'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.65%  '

# Row 29
$ws.Range("B29").Value = 'InjectiveProtocol'
$ws.Range("C29").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '42.44'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +11.04%  '

# Row 30
$ws.Range("E30").Value = '  -1.97%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '2.27'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.06%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '173.40'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.76%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0936'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +5.15%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '21.62'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.83%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '5.71'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.52%  '

# Row 36
$ws.Range("E36").Value = '  +2.12%  '

# Row 38
$ws.Range("E38").Value = '  +3.21%  '

# Row 39
$ws.Range("E39").Value = '  -1.55%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.86'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +9.18%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '14.24'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +17.06%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '74.68'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +15.71%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.43'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +4.52%  '

# Row 44
$ws.Range("E44").Value = '  +1.41%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '6.37'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +22.40%  '

# Row 46
$ws.Range("E46").Value = '  +0.02%  '

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.40'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +1.28%  '

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '8.76'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.12%  '

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '102.90'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.21%  '

# Row 50
$ws.Range("E50").Value = '  +3.62%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.1000'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.91%  '
